$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: header strings (B2:H2) to seed shared-string order ---
$ws.Range("B2").Value = "anchor score"
$ws.Range("C2").Value = "type occurences"
$ws.Range("D2").Value = "total occurences"
$ws.Range("E2").Value = "+%"
$ws.Range("F2").Value = "-%"
$ws.Range("G2").Value = "both"
$ws.Range("H2").Value = "normal"

# --- Step 2: column A labels, row 2 (name) then 3-35 top to bottom ---
$ws.Range("A2").Value = "name"
$ws.Range("A3").Value = "poorly"
$ws.Range("A4").Value = "disappointing"
$ws.Range("A5").Value = "poor"
$ws.Range("A6").Value = "junk"
$ws.Range("A7").Value = "disappointed"
$ws.Range("A8").Value = "broke"
$ws.Range("A9").Value = "however"
$ws.Range("A10").Value = "waste"
$ws.Range("A11").Value = "smaller"
$ws.Range("A12").Value = "small"
$ws.Range("A13").Value = "broken"
$ws.Range("A14").Value = "apart"
$ws.Range("A15").Value = "plastic"
$ws.Range("A16").Value = "cheap"
$ws.Range("A17").Value = "difficult"
$ws.Range("A18").Value = "ok"
$ws.Range("A19").Value = "thought"
$ws.Range("A20").Value = "though"
$ws.Range("A21").Value = "size"
$ws.Range("A22").Value = "hard"
$ws.Range("A23").Value = "item"
$ws.Range("A24").Value = "would"
$ws.Range("A25").Value = "money"
$ws.Range("A26").Value = "better"
$ws.Range("A27").Value = "price"
$ws.Range("A28").Value = "work"
$ws.Range("A29").Value = "product"
$ws.Range("A30").Value = "use"
$ws.Range("A31").Value = "buy"
$ws.Range("A32").Value = "little"
$ws.Range("A33").Value = "like"
$ws.Range("A34").Value = "one"
$ws.Range("A35").Value = "toy"

# --- Step 3: A1 (negative) ---
$ws.Range("A1").Value = "negative"

# --- Step 4: column J labels, row 2 (name) then 3-24 top to bottom ---
$ws.Range("J2").Value = "name"
$ws.Range("J3").Value = "wonderful"
$ws.Range("J4").Value = "awesome"
$ws.Range("J5").Value = "favorite"
$ws.Range("J6").Value = "excellent"
$ws.Range("J7").Value = "classic"
$ws.Range("J8").Value = "thank"
$ws.Range("J9").Value = "love"
$ws.Range("J10").Value = "loves"
$ws.Range("J11").Value = "great"
$ws.Range("J12").Value = "loved"
$ws.Range("J13").Value = "friends"
$ws.Range("J14").Value = "perfect"
$ws.Range("J15").Value = "best"
$ws.Range("J16").Value = "learn"
$ws.Range("J17").Value = "happy"
$ws.Range("J18").Value = "enjoy"
$ws.Range("J19").Value = "christmas"
$ws.Range("J20").Value = "fun"
$ws.Range("J21").Value = "game"
$ws.Range("J22").Value = "family"
$ws.Range("J23").Value = "easy"
$ws.Range("J24").Value = "play"

# --- Step 5: J1 (positive) ---
$ws.Range("J1").Value = "positive"
# --- Step 6: row 2 mirrored headers K2:Q2 (reuse existing shared strings) ---
$ws.Range("J2").Value = "name"
$ws.Range("K2").Value = "anchor score"
$ws.Range("L2").Value = "type occurences"
$ws.Range("M2").Value = "total occurences"
$ws.Range("N2").Value = "+%"
$ws.Range("O2").Value = "-%"
$ws.Range("P2").Value = "both"
$ws.Range("Q2").Value = "normal"

# --- Step 7: numeric/boolean data cells, columns B-H for rows 3-35 ---
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 46
$ws.Range("D3").Value = 46
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 0
$ws.Range("B4").Value = 0.7727272727272727
$ws.Range("C4").Value = 34
$ws.Range("D4").Value = 34
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 10
$ws.Range("B5").Value = 0.7183098591549296
$ws.Range("C5").Value = 51
$ws.Range("D5").Value = 51
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 20
$ws.Range("B6").Value = 0.6909090909090909
$ws.Range("C6").Value = 38
$ws.Range("D6").Value = 38
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 17
$ws.Range("B7").Value = 0.6774193548387096
$ws.Range("C7").Value = 126
$ws.Range("D7").Value = 126
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 60
$ws.Range("B8").Value = 0.6504854368932039
$ws.Range("C8").Value = 134
$ws.Range("D8").Value = 134
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 72
$ws.Range("B9").Value = 0.640625
$ws.Range("C9").Value = 41
$ws.Range("D9").Value = 41
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = 23
$ws.Range("B10").Value = 0.6283783783783784
$ws.Range("C10").Value = 93
$ws.Range("D10").Value = 93
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = $false
$ws.Range("H10").Value = 55
$ws.Range("B11").Value = 0.5882352941176471
$ws.Range("C11").Value = 70
$ws.Range("D11").Value = 70
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = $false
$ws.Range("H11").Value = 49
$ws.Range("B12").Value = 0.5072463768115942
$ws.Range("C12").Value = 175
$ws.Range("D12").Value = 175
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = $false
$ws.Range("H12").Value = 170
$ws.Range("B13").Value = 0.4337349397590362
$ws.Range("C13").Value = 36
$ws.Range("D13").Value = 36
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = $false
$ws.Range("H13").Value = 47
$ws.Range("B14").Value = 0.4210526315789473
$ws.Range("C14").Value = 40
$ws.Range("D14").Value = 40
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = $false
$ws.Range("H14").Value = 55
$ws.Range("B15").Value = 0.4173228346456693
$ws.Range("C15").Value = 53
$ws.Range("D15").Value = 53
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = $false
$ws.Range("H15").Value = 74
$ws.Range("B16").Value = 0.3601895734597156
$ws.Range("C16").Value = 76
$ws.Range("D16").Value = 76
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = $false
$ws.Range("H16").Value = 135
$ws.Range("B17").Value = 0.3258426966292135
$ws.Range("C17").Value = 29
$ws.Range("D17").Value = 29
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = $false
$ws.Range("H17").Value = 60
$ws.Range("B18").Value = 0.3125
$ws.Range("C18").Value = 40
$ws.Range("D18").Value = 40
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = $false
$ws.Range("H18").Value = 88
$ws.Range("B19").Value = 0.3069306930693069
$ws.Range("C19").Value = 62
$ws.Range("D19").Value = 62
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = $false
$ws.Range("H19").Value = 140
$ws.Range("B20").Value = 0.2478632478632479
$ws.Range("C20").Value = 29
$ws.Range("D20").Value = 29
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = $false
$ws.Range("H20").Value = 88
$ws.Range("B21").Value = 0.2319587628865979
$ws.Range("C21").Value = 45
$ws.Range("D21").Value = 45
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = $false
$ws.Range("H21").Value = 149
$ws.Range("B22").Value = 0.21
$ws.Range("C22").Value = 42
$ws.Range("D22").Value = 42
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = $false
$ws.Range("H22").Value = 158
$ws.Range("B23").Value = 0.1992753623188406
$ws.Range("C23").Value = 55
$ws.Range("D23").Value = 55
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = $false
$ws.Range("H23").Value = 221
$ws.Range("B24").Value = 0.187221396731055
$ws.Range("C24").Value = 126
$ws.Range("D24").Value = 127
$ws.Range("E24").Value = 0.01
$ws.Range("F24").Value = 0.99
$ws.Range("G24").Value = $true
$ws.Range("H24").Value = 547
$ws.Range("B25").Value = 0.1867088607594937
$ws.Range("C25").Value = 59
$ws.Range("D25").Value = 59
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = $false
$ws.Range("H25").Value = 257
$ws.Range("B26").Value = 0.1635514018691589
$ws.Range("C26").Value = 35
$ws.Range("D26").Value = 35
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = $false
$ws.Range("H26").Value = 179
$ws.Range("B27").Value = 0.1613832853025937
$ws.Range("C27").Value = 56
$ws.Range("D27").Value = 57
$ws.Range("E27").Value = 0.02
$ws.Range("F27").Value = 0.98
$ws.Range("G27").Value = $true
$ws.Range("H27").Value = 291
$ws.Range("B28").Value = 0.1587301587301587
$ws.Range("C28").Value = 50
$ws.Range("D28").Value = 51
$ws.Range("E28").Value = 0.02
$ws.Range("F28").Value = 0.98
$ws.Range("G28").Value = $true
$ws.Range("H28").Value = 265
$ws.Range("B29").Value = 0.1431718061674009
$ws.Range("C29").Value = 65
$ws.Range("D29").Value = 65
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = $false
$ws.Range("H29").Value = 389
$ws.Range("B30").Value = 0.09315068493150686
$ws.Range("C30").Value = 34
$ws.Range("D30").Value = 34
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = $false
$ws.Range("H30").Value = 331
$ws.Range("B31").Value = 0.09295774647887324
$ws.Range("C31").Value = 33
$ws.Range("D31").Value = 33
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = $false
$ws.Range("H31").Value = 322
$ws.Range("B32").Value = 0.08928571428571429
$ws.Range("C32").Value = 40
$ws.Range("D32").Value = 41
$ws.Range("E32").Value = 0.02
$ws.Range("F32").Value = 0.98
$ws.Range("G32").Value = $true
$ws.Range("H32").Value = 408
$ws.Range("B33").Value = 0.07260726072607261
$ws.Range("C33").Value = 44
$ws.Range("D33").Value = 46
$ws.Range("E33").Value = 0.04
$ws.Range("F33").Value = 0.96
$ws.Range("G33").Value = $true
$ws.Range("H33").Value = 562
$ws.Range("B34").Value = 0.05979643765903308
$ws.Range("C34").Value = 47
$ws.Range("D34").Value = 55
$ws.Range("E34").Value = 0.15
$ws.Range("F34").Value = 0.85
$ws.Range("G34").Value = $true
$ws.Range("H34").Value = 739
$ws.Range("B35").Value = 0.05053598774885146
$ws.Range("C35").Value = 33
$ws.Range("D35").Value = 35
$ws.Range("E35").Value = 0.06
$ws.Range("F35").Value = 0.9399999999999999
$ws.Range("G35").Value = $true
$ws.Range("H35").Value = 620

# --- Step 8: numeric/boolean data cells, columns K-Q for rows 3-24 ---
$ws.Range("K3").Value = 0.8928571428571429
$ws.Range("L3").Value = 50
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 6
$ws.Range("K4").Value = 0.8615384615384616
$ws.Range("L4").Value = 56
$ws.Range("M4").Value = 56
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 9
$ws.Range("K5").Value = 0.8172043010752689
$ws.Range("L5").Value = 76
$ws.Range("M5").Value = 76
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 17
$ws.Range("K6").Value = 0.796875
$ws.Range("L6").Value = 51
$ws.Range("M6").Value = 51
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 13
$ws.Range("K7").Value = 0.6981132075471698
$ws.Range("L7").Value = 37
$ws.Range("M7").Value = 37
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 16
$ws.Range("K8").Value = 0.5797101449275363
$ws.Range("L8").Value = 40
$ws.Range("M8").Value = 40
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 29
$ws.Range("K9").Value = 0.5796269727403156
$ws.Range("L9").Value = 404
$ws.Range("M9").Value = 404
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 293
$ws.Range("K10").Value = 0.5456431535269709
$ws.Range("L10").Value = 263
$ws.Range("M10").Value = 263
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 219
$ws.Range("K11").Value = 0.4716981132075472
$ws.Range("L11").Value = 575
$ws.Range("M11").Value = 576
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $true
$ws.Range("Q11").Value = 644
$ws.Range("K12").Value = 0.3730886850152905
$ws.Range("L12").Value = 122
$ws.Range("M12").Value = 122
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 205
$ws.Range("K13").Value = 0.3650793650793651
$ws.Range("L13").Value = 69
$ws.Range("M13").Value = 69
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 120
$ws.Range("K14").Value = 0.3614457831325301
$ws.Range("L14").Value = 60
$ws.Range("M14").Value = 60
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 106
$ws.Range("K15").Value = 0.35
$ws.Range("L15").Value = 42
$ws.Range("M15").Value = 42
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 78
$ws.Range("K16").Value = 0.2734375
$ws.Range("L16").Value = 35
$ws.Range("M16").Value = 35
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 93
$ws.Range("K17").Value = 0.2097902097902098
$ws.Range("L17").Value = 30
$ws.Range("M17").Value = 30
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 113
$ws.Range("K18").Value = 0.1989247311827957
$ws.Range("L18").Value = 37
$ws.Range("M18").Value = 37
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 149
$ws.Range("K19").Value = 0.1967871485943775
$ws.Range("L19").Value = 49
$ws.Range("M19").Value = 49
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 200
$ws.Range("K20").Value = 0.1763157894736842
$ws.Range("L20").Value = 201
$ws.Range("M20").Value = 202
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $true
$ws.Range("Q20").Value = 939
$ws.Range("K21").Value = 0.1027308192457737
$ws.Range("L21").Value = 158
$ws.Range("M21").Value = 161
$ws.Range("N21").Value = 0.98
$ws.Range("O21").Value = 0.02000000000000002
$ws.Range("P21").Value = $true
$ws.Range("Q21").Value = 1380
$ws.Range("K22").Value = 0.0947075208913649
$ws.Range("L22").Value = 34
$ws.Range("M22").Value = 34
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = $false
$ws.Range("Q22").Value = 325
$ws.Range("K23").Value = 0.09358288770053476
$ws.Range("L23").Value = 35
$ws.Range("M23").Value = 35
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = $false
$ws.Range("Q23").Value = 339
$ws.Range("K24").Value = 0.04806408544726302
$ws.Range("L24").Value = 36
$ws.Range("M24").Value = 39
$ws.Range("N24").Value = 0.92
$ws.Range("O24").Value = 0.07999999999999996
$ws.Range("P24").Value = $true
$ws.Range("Q24").Value = 713

# --- Step 9: remove now-obsolete row 36 ---
$ws.Rows(36).Delete()
